# ------------------------------------------------------------------
# MapleStory Monthly Honorable Rock Record — add "April 2020" sheet
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ----- 1. Normalise the saved selection on each existing month sheet -----
# (each was previously a two-piece sqref left over from an old selection;
#  collapse it back down to the single active cell, sheet by sheet)
$selFixes = @(
  @("August 2019",    "D4"),
  @("September 2019", "D4"),
  @("October 2019",   "D4"),
  @("November 2019",  "D5"),
  @("December 2019",  "D5"),
  @("January 2020",   "D6"),
  @("February 2020",  "E2"),
  @("March 2020",     "E2")
)

foreach ($fix in $selFixes) {
    $sheetName = $fix[0]
    $cellRef = $fix[1]
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Activate() | Out-Null
    $ws.Range($cellRef).Select() | Out-Null
}

# ----- 2. Create the new "April 2020" sheet from the "March 2020" template -----
$template = $wb.Worksheets.Item("March 2020")
$template.Copy([System.Reflection.Missing]::Value, $template)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "April 2020"

# Update the embedded title text (row 2, col E) from "March 2020" to "April 2020"
$ws.Range("E2").Formula = '="    """&"April 2020"&""""&":"'

# ----- 3. Replace the guild/contribution rows (4-53) with the April 2020 results -----
# Column B (rank 1..50) is already correct after the sheet copy.
$data = @(
  @('Savages', 261689862),
  @('Smile', 259630009),
  @('Eternal', 255017557),
  @('Bounce', 230945090),
  @('Elite', 212697617),
  @('Epic', 149588101),
  @('Spring', 148563900),
  @('Sunset', 142194188),
  @('Beaters', 132902319),
  @('Downtime', 126969925),
  @('lolicafe', 108673933),
  @('Imperium', 104985564),
  @('Remorse', 98790172),
  @('Cleanse', 91543234),
  @('RainSong', 89789290),
  @('Gintama', 88678345),
  @('Revive', 88630272),
  @('Undertale', 88410162),
  @('Maha', 87440842),
  @('Lithe', 85509331),
  @('Broke', 83738795),
  @('Sora', 80147968),
  @('Erda', 79239156),
  @('Rising', 67788196),
  @('Sugar', 67402226),
  @('Earnest', 62525517),
  @('Ravers', 60420653),
  @('Oceania', 59452636),
  @('Aloe', 58638996),
  @('Fabled', 57269706),
  @('Howl', 56903428),
  @('Mystical', 53603764),
  @('CyberThreat', 53262304),
  @('Atelier', 51707222),
  @('RainDrop', 51698184),
  @('Kingdom', 48960377),
  @('Tama', 48347462),
  @('Exorcist', 47182813),
  @('Fandom', 45631841),
  @('Nutsy', 45518787),
  @('Weibo', 43698709),
  @('Path', 42432322),
  @('Miao', 42334268),
  @('Bubbles', 42167767),
  @('Reboot', 40987120),
  @('Artifacts', 40697236),
  @('chigga', 40204303),
  @('Comity', 39916833),
  @('Gentle', 39185113),
  @('SlimeTree', 38860946)
)

$r = 4
foreach ($item in $data) {
    $ws.Range("C$r").Value = $item[0]
    $ws.Range("D$r").Value = $item[1]
    $r = $r + 1
}

# ----- 4. Make the new sheet the active tab with its own saved selection -----
$ws.Activate() | Out-Null
$ws.Range("E55").Select() | Out-Null
